$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "91.414.08"
$ws.Range("E2").Value = "  +0.92%  "

# Row 3
$ws.Range("D3").Value = "3.168.20"
$ws.Range("E3").Value = "  +1.88%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.66"
$ws.Range("E5").Value = "  +2.34%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "621.78"
$ws.Range("E6").Value = "  -0.25%  "

# Row 7
$ws.Range("E7").Value = "  +5.43%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.374"
$ws.Range("E8").Value = "  +1.86%  "

# Row 9
$ws.Range("E9").Value = "  -0.17%  "

# Row 10
$ws.Range("D10").Value = "3.164.74"
$ws.Range("E10").Value = "  +1.81%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.750"
$ws.Range("E11").Value = "  +1.18%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.206"
$ws.Range("E12").Value = "  +4.13%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000249"
$ws.Range("E13").Value = "  -1.11%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.54"
$ws.Range("E14").Value = "  -1.27%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.54"
$ws.Range("E15").Value = "  +1.05%  "

# Row 16
$ws.Range("D16").Value = "91.223.85"
$ws.Range("E16").Value = "  +0.99%  "

# Row 17
$ws.Range("D17").Value = "3.753.17"
$ws.Range("E17").Value = "  +2.39%  "

# Row 18
$ws.Range("D18").Value = "3.169.32"
$ws.Range("E18").Value = "  +2.82%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.74"
$ws.Range("E19").Value = "  -4.43%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.34"
$ws.Range("E20").Value = "  +9.13%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.91"
$ws.Range("E21").Value = "  +5.97%  "

# Row 22
$ws.Range("E22").Value = "  -5.87%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "443.87"
$ws.Range("E23").Value = "  +1.07%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.22"
$ws.Range("E24").Value = "  +2.86%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.10"
$ws.Range("E25").Value = "  +3.02%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "88.96"
$ws.Range("E26").Value = "  +0.24%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.17"
$ws.Range("E27").Value = "  +0.37%  "

# Row 28
$ws.Range("D28").Value = "3.349.27"
$ws.Range("E28").Value = "  +2.51%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.19%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.129"
$ws.Range("E30").Value = "  +45.41%  "

# Row 31
$ws.Range("E31").Value = "  +6.20%  "

# Row 32
$ws.Range("E32").Value = "  +17.06%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "9.40"
$ws.Range("E33").Value = "  +0.30%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.168"
$ws.Range("E34").Value = "  +9.94%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.75"
$ws.Range("E35").Value = "  +6.29%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.50"
$ws.Range("E36").Value = "  +2.03%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.905"
$ws.Range("E37").Value = "  -9.65%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "512.70"
$ws.Range("E38").Value = "  +1.53%  "

# Row 39
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.37"
$ws.Range("E39").Value = "  +6.19%  "

# Row 40
$ws.Range("B40").Value = "PancakeSwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.94"
$ws.Range("E40").Value = "  +1.54%  "

# Row 41
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.454"
$ws.Range("E41").Value = "  +11.55%  "

# Row 42
$ws.Range("B42").Value = "MantraDAO"
$ws.Range("C42").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.85"
$ws.Range("E42").Value = "  +11.98%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.46"
$ws.Range("E43").Value = "  -10.73%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.14"
$ws.Range("E44").Value = "  -0.15%  "

# Row 45
$ws.Range("E45").Value = "  +0.01%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.723"
$ws.Range("E46").Value = "  +4.31%  "

# Row 47
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.93"
$ws.Range("E47").Value = "  +1.09%  "

# Row 48
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "156.67"
$ws.Range("E48").Value = "  +3.72%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.40"
$ws.Range("E49").Value = "  +4.19%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.47"
$ws.Range("E50").Value = "  +0.32%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.06"
$ws.Range("E51").Value = "  -1.29%  "
